$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = -12.294
$ws.Range("A3").Value = -21.711
$ws.Range("D3").Value = -7.827000000000001
$ws.Range("E6").Value = 16.844
$ws.Range("D12").Value = -7.347
$ws.Range("A14").Value = -21.668
$ws.Range("E19").Value = 16.538
$ws.Range("A21").Value = -20.21
$ws.Range("A23").Value = -20.142
$ws.Range("D24").Value = -7.527000000000001
$ws.Range("E24").Value = 17
$ws.Range("A25").Value = -22.017
$ws.Range("C25").Value = -13.169
$ws.Range("D25").Value = -8.717000000000001
$ws.Range("A26").Value = -21.611
$ws.Range("C27").Value = -13.163
$ws.Range("A29").Value = -21.284
$ws.Range("E30").Value = 16.348
$ws.Range("C31").Value = -12.532
$ws.Range("E31").Value = 16.549
$ws.Range("E33").Value = 17.573
$ws.Range("C39").Value = -12.847
$ws.Range("E42").Value = 16.727
$ws.Range("C48").Value = -11.598
$ws.Range("D50").Value = -8.270000000000001
$ws.Range("C51").Value = -11.564
$ws.Range("C52").Value = -11.606
$ws.Range("A53").Value = -22.128
$ws.Range("D53").Value = -7.320000000000002
$ws.Range("C55").Value = -13.412
$ws.Range("E55").Value = 16.349
$ws.Range("C56").Value = -12.753
$ws.Range("A57").Value = -21.891
$ws.Range("C57").Value = -12.926
$ws.Range("D57").Value = -8.606
$ws.Range("E58").Value = 17.007
$ws.Range("A59").Value = -22.358
$ws.Range("D61").Value = -7.708
$ws.Range("D63").Value = -7.811
$ws.Range("E65").Value = 17.151
$ws.Range("A69").Value = -21.68
$ws.Range("D70").Value = -7.356999999999999
$ws.Range("E70").Value = 17.612
$ws.Range("C73").Value = -12.763
$ws.Range("E75").Value = 16.76
$ws.Range("A79").Value = -20.849
$ws.Range("A83").Value = -21.938
$ws.Range("E83").Value = 16.737
$ws.Range("D86").Value = -8.322000000000001
$ws.Range("E86").Value = 16.323
$ws.Range("C89").Value = -12.39
$ws.Range("C90").Value = -12.942
$ws.Range("A91").Value = -21.509
$ws.Range("C92").Value = -11.823
$ws.Range("A93").Value = -21.318
$ws.Range("E96").Value = 16.454
$ws.Range("E97").Value = 16.846
$ws.Range("D98").Value = -8.361000000000001
$ws.Range("D100").Value = -8.409000000000001
$ws.Range("D102").Value = -7.794000000000001
